$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update Id, Taxonsorteringsordning, Rodlistade, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Starttid, Sluttid
$ws.Range("A2").Value = 111708099
$ws.Range("B2").Value = 90660
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 4362
$ws.Range("F2").Value = "Blå taggsvamp"
$ws.Range("G2").Value = "Hydnellum caeruleum"
$ws.Range("H2").Value = "(Hornem.) P.Karst."
$ws.Range("Z2").Value = "14:16"
$ws.Range("AB2").Value = "14:16"

# Row 7: update Id, Taxonsorteringsordning, Rodlistade, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Starttid, Sluttid
$ws.Range("A7").Value = 111708920
$ws.Range("B7").Value = 90666
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 4364
$ws.Range("F7").Value = "Dropptaggsvamp"
$ws.Range("G7").Value = "Hydnellum ferrugineum"
$ws.Range("H7").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Z7").Value = "13:53"
$ws.Range("AB7").Value = "13:53"
